$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column C ("Förändrad") changes from 45184 to 45186 for every existing data row (2-407)
$ws.Range("C2:C407").Value = 45186

# 2. Add the display-text second argument to every HYPERLINK() formula so the
#    link shows the case "Beteckning" (column A) instead of the raw URL.
#    Rows 2-19 have S,T,V,W,X,Y (and row 10 additionally has U).
#    Rows 257 and 292 only have U,V,W,X,Y.
$hyperlinkRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,257,292)
$hyperlinkCols = @("S","T","U","V","W","X","Y")

foreach ($r in $hyperlinkRows) {
    $name = $ws.Range("A$r").Value()
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula()
        if ($f -and $f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
            $cell.Formula = $newFormula
        }
    }
}

# 3. Row 407 gains an explicit row-height attribute (ht="15" customHeight="1")
$ws.Rows.Item(407).RowHeight = 15

# 4. Append a brand-new row 408 for case "A 43455-2023"
$ws.Range("A408").Value = "A 43455-2023"
$ws.Range("B408").Value = 45184
$ws.Range("C408").Value = 45186
$ws.Range("D408").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E408").Value = "YDRE"
$ws.Range("G408").Value = 1.8
$ws.Range("H408").Value = 0
$ws.Range("I408").Value = 0
$ws.Range("J408").Value = 0
$ws.Range("K408").Value = 0
$ws.Range("L408").Value = 0
$ws.Range("M408").Value = 0
$ws.Range("N408").Value = 0
$ws.Range("O408").Value = 0
$ws.Range("P408").Value = 0
$ws.Range("Q408").Value = 0
$ws.Range("R408").Value = ""

# Match the date-number formatting and wrap-text style used by the other rows
$ws.Range("B408").NumberFormat = $ws.Range("B407").NumberFormat()
$ws.Range("C408").NumberFormat = $ws.Range("C407").NumberFormat()
$ws.Range("R408").WrapText = $true
